# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-price update to the Kujata_Profits workbook.
# For each affected leve row, updates price/profit columns (H..N) to the refreshed
# market-board snapshot, clearing any column that no longer has a value in the new data.

$wb = $excel.ActiveWorkbook

# ALC!row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2498.617
$ws.Range("I98").Value = 2613.0952
$ws.Range("J98").Value = 1537
$ws.Range("K98").Value = 2613.0952
$ws.Range("L98").Value = 1537
$ws.Range("M98").Value = -1115.0952
$ws.Range("N98").Value = -4533

# ALC!row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2498.617
$ws.Range("I122").Value = 2613.0952
$ws.Range("J122").Value = 1537
$ws.Range("K122").Value = 7839.285600000001
$ws.Range("L122").Value = 4611
$ws.Range("M122").Value = -5389.285600000001
$ws.Range("N122").Value = -9511

# ALC!row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7944320
$ws.Range("I132").Value = 15880564
$ws.Range("J132").Value = 8075.7617
$ws.Range("K132").Value = 47641692
$ws.Range("L132").Value = 24227.2851
$ws.Range("M132").Value = -47639162
$ws.Range("N132").Value = -29287.2851

# ARM!row 10
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6372.64
$ws.Range("I32").Value = 5368.591
$ws.Range("K32").Value = 5368.591
$ws.Range("M32").Value = -5081.591

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3172.3572
$ws.Range("I132").Value = 2450
$ws.Range("J132").Value = 4978.25
$ws.Range("K132").Value = 7350
$ws.Range("L132").Value = 14934.75
$ws.Range("M132").Value = -4820
$ws.Range("N132").Value = -19994.75

# BSM!row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 487.5
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -377
$ws.Range("N22").Value = -646

# BSM!row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6370.095
$ws.Range("I134").Value = 1079.5625
$ws.Range("K134").Value = 3238.6875
$ws.Range("M134").Value = -703.6875

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1298.7142
$ws.Range("I31").Value = 1267.7142
$ws.Range("J31").Value = 1360.7142
$ws.Range("K31").Value = 1267.7142
$ws.Range("L31").Value = 1360.7142
$ws.Range("M31").Value = -972.7141999999999
$ws.Range("N31").Value = -1950.7142

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1298.7142
$ws.Range("I34").Value = 1267.7142
$ws.Range("J34").Value = 1360.7142
$ws.Range("K34").Value = 1267.7142
$ws.Range("L34").Value = 1360.7142
$ws.Range("M34").Value = -1065.7142
$ws.Range("N34").Value = -1764.7142

# CRP!row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 943
$ws.Range("I122").Value = 900.2727
$ws.Range("J122").Value = 1099.6666
$ws.Range("K122").Value = 2700.8181
$ws.Range("L122").Value = 3298.9998
$ws.Range("M122").Value = -250.8181
$ws.Range("N122").Value = -8198.9998

# CRP!row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1788.4572
$ws.Range("I132").Value = 1464.7142
$ws.Range("K132").Value = 4394.142599999999
$ws.Range("M132").Value = -1864.142599999999

# CRP!row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 16130934
$ws.Range("I134").Value = 1938.9584
$ws.Range("J134").Value = 71430344
$ws.Range("K134").Value = 5816.8752
$ws.Range("L134").Value = 214291032
$ws.Range("M134").Value = -3281.8752
$ws.Range("N134").Value = -214296102

# CRP!row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 37950
$ws.Range("J140").Value = 37950
$ws.Range("L140").Value = 37950
$ws.Range("N140").Value = -48310

# CRP!row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 2999900
$ws.Range("J141").Value = 2999900
$ws.Range("L141").Value = 2999900
$ws.Range("N141").Value = -3010260

# CUL!row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 12822986
$ws.Range("J131").Value = 2719.0286
$ws.Range("L131").Value = 8157.085800000001
$ws.Range("N131").Value = -18237.0858

# GSM!row 19
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 970
$ws.Range("I19").Value = 970
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 970
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -682
$ws.Range("N19").ClearContents()

# GSM!row 21
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 2503200
$ws.Range("J21").Value = 6400
$ws.Range("L21").Value = 6400
$ws.Range("N21").Value = -6746

# GSM!row 30
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 2503200
$ws.Range("J30").Value = 6400
$ws.Range("L30").Value = 6400
$ws.Range("N30").Value = -6610

# LTW!row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2222.4443
$ws.Range("I7").Value = 1800.6
$ws.Range("J7").Value = 2749.75
$ws.Range("K7").Value = 1800.6
$ws.Range("L7").Value = 2749.75
$ws.Range("M7").Value = -1688.6
$ws.Range("N7").Value = -2973.75

# LTW!row 10
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = 30
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 30
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 110
$ws.Range("N10").ClearContents()

# LTW!row 25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 2833.3333
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 2833.3333
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 2833.3333
$ws.Range("N25").Value = -3293.3333
$ws.Range("M25").ClearContents()

# LTW!row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1038.875
$ws.Range("I55").Value = 901.5714
$ws.Range("J55").Value = 2000
$ws.Range("K55").Value = 901.5714
$ws.Range("L55").Value = 2000
$ws.Range("M55").Value = -728.5714
$ws.Range("N55").Value = -2346

# LTW!row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1790.8125
$ws.Range("I82").Value = 1641.3636
$ws.Range("J82").Value = 2119.6
$ws.Range("K82").Value = 1641.3636
$ws.Range("L82").Value = 2119.6
$ws.Range("M82").Value = -1280.3636
$ws.Range("N82").Value = -2841.6

# LTW!row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1790.8125
$ws.Range("I85").Value = 1641.3636
$ws.Range("J85").Value = 2119.6
$ws.Range("K85").Value = 1641.3636
$ws.Range("L85").Value = 2119.6
$ws.Range("M85").Value = -393.3635999999999
$ws.Range("N85").Value = -4615.6

# LTW!row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2222.4443
$ws.Range("I126").Value = 1800.6
$ws.Range("J126").Value = 2749.75
$ws.Range("K126").Value = 5401.799999999999
$ws.Range("L126").Value = 8249.25
$ws.Range("M126").Value = -2931.799999999999
$ws.Range("N126").Value = -13189.25

# WVR!row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2088.25
$ws.Range("I132").Value = 1831
$ws.Range("J132").Value = 3374.5
$ws.Range("K132").Value = 5493
$ws.Range("L132").Value = 10123.5
$ws.Range("M132").Value = -2963
$ws.Range("N132").Value = -15183.5

